$d = $word.ActiveDocument

# Find the end of the target paragraph text and position the range there.
$rng = $d.Content
$found = $rng.Find.Execute("and jumper cables, powered by a power bank on the right strap and controlled by the Arduino on the left strap below the breadboard.", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)

if ($found) {
    # Collapse to the end of the found text
    $rng.Collapse(0)

    # Insert an empty paragraph, then the Sprint 4 Review heading paragraph.
    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(4, 1) | Out-Null

    $rng.InsertParagraphAfter()
    $rng.Collapse(0)
    $rng.Move(4, 1) | Out-Null

    $rng.Text = "Sprint 4 Review"
    $rng.Font.Size = 14
    $rng.Font.Underline = 1
    $rng.ParagraphFormat.Alignment = 1
}

Write-Output "done"
